$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.228.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.84%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.158.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -8.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.76%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.28%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.39%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.154.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.45%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.98%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.393'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.683.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.99%  '

# Row 14
$ws.Range("E14").Value = '  +0.48%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -9.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.091.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.96%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000163'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.03%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.142.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.40%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.08%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.45%  '

# Row 23
$ws.Range("E23").Value = '  -0.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.54%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.499'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000116'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -10.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.59'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.175'
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("E30").Value = '  -0.13%  '

# Row 31
$ws.Range("E31").Value = '  -5.81%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.79%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.54%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.99%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.821'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.97%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.641.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.90%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.85%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.99%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.92%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0650'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.00%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.79%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '315.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.56%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0270'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.30%  '

# Row 50
$ws.Range("E50").Value = '  -4.62%  '

# Row 51
$ws.Range("E51").Value = '  -0.13%  '
